$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with revised AgTests/AgPosit figures ---
$ws.Range("F523").Value = 10191

$ws.Range("F530").Value = 12643
$ws.Range("G530").Value = 42

$ws.Range("F531").Value = 9079

$ws.Range("F532").Value = 10157

$ws.Range("F533").Value = 11741

$ws.Range("F534").Value = 16549
$ws.Range("G534").Value = 51

$ws.Range("F535").Value = 9996

$ws.Range("F536").Value = 7859

$ws.Range("F537").Value = 13405

$ws.Range("F538").Value = 10735

$ws.Range("F539").Value = 10067
$ws.Range("G539").Value = 44

# --- Append new rows 540-543 with data through 30.08.2021 ---
$ws.Range("A540").Value = 44434
$ws.Range("B540").Value = 394535
$ws.Range("C540").Value = 5017
$ws.Range("D540").Value = 89
$ws.Range("E540").Value = 12548
$ws.Range("F540").Value = 11697
$ws.Range("G540").Value = 61

$ws.Range("A541").Value = 44435
$ws.Range("B541").Value = 394657
$ws.Range("C541").Value = 6542
$ws.Range("D541").Value = 122
$ws.Range("E541").Value = 12548
$ws.Range("F541").Value = 14450
$ws.Range("G541").Value = 58

$ws.Range("A542").Value = 44436
$ws.Range("B542").Value = 394742
$ws.Range("C542").Value = 2832
$ws.Range("D542").Value = 85
$ws.Range("E542").Value = 12548
$ws.Range("F542").Value = 8622
$ws.Range("G542").Value = 44

$ws.Range("A543").Value = 44437
$ws.Range("B543").Value = 394791
$ws.Range("C543").Value = 1485
$ws.Range("D543").Value = 49
$ws.Range("E543").Value = 12548
$ws.Range("F543").Value = 3548
$ws.Range("G543").Value = 25

# Ensure new date cells (column A) use the yyyy-mm-dd date format,
# matching the rest of the column.
$ws.Range("A540:A543").NumberFormat = "yyyy-mm-dd"
